$wb = $excel.ActiveWorkbook

# "AddCustomerTest" sheet holds the test data for the add-customer scenario.
# Its runmode flag (row 2, column A) is being switched on from "N" to "Y".
$ws = $wb.Worksheets.Item("AddCustomerTest")
$ws.Range("A2").Value = "Y"

# Move the active selection to A4, as reflected in the saved view state.
$ws.Range("A4").Select()
